# The commit swaps slides 17 and 18 (their SlideIDs 723/724 trade places
# while the rest of the deck keeps its order), so slide 17's content
# ("UML Diagram Types" overview) becomes slide 18 and slide 18's content
# ("The Component Diagram" detail) becomes slide 17.
#
# PowerPoint's Slide.MoveTo(newIndex) is the native COM operation for
# reordering slides and produces exactly that end state (SlideID/content
# pairing swapped between positions 17 and 18, every other slide untouched).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(17)
$s.MoveTo(18)
